$d = $word.ActiveDocument

$replacements = @(
    @("959÷9=106, 5", "261÷3=87, 0"),
    @("650÷8=81, 2", "618÷6=103, 0"),
    @("748÷4=187, 0", "811÷7=115, 6"),
    @("939÷6=156, 3", "436÷6=72, 4"),
    @("752÷5=150, 2", "796÷6=132, 4"),
    @("461÷6=76, 5", "746÷9=82, 8"),
    @("180÷3=60, 0", "113÷4=28, 1"),
    @("155÷7=22, 1", "351÷8=43, 7"),
    @("981÷6=163, 3", "284÷4=71, 0"),
    @("914÷2=457, 0", "904÷9=100, 4"),
    @("965÷8=120, 5", "328÷8=41, 0"),
    @("346÷4=86, 2", "199÷5=39, 4"),
    @("580÷9=64, 4", "749÷9=83, 2"),
    @("609÷6=101, 3", "301÷3=100, 1"),
    @("148÷4=37, 0", "160÷2=80, 0"),
    @("514÷2=257, 0", "851÷7=121, 4"),
    @("627÷8=78, 3", "511÷7=73, 0"),
    @("686÷9=76, 2", "276÷5=55, 1"),
    @("338÷8=42, 2", "706÷6=117, 4"),
    @("690÷4=172, 2", "955÷6=159, 1"),
    @("761÷3=253, 2", "274÷3=91, 1"),
    @("718÷8=89, 6", "631÷8=78, 7"),
    @("625÷5=125, 0", "476÷2=238, 0"),
    @("116÷3=38, 2", "460÷6=76, 4"),
    @("119÷2=59, 1", "814÷2=407, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
